$wb = $excel.ActiveWorkbook

# --- "Sheet4" tab (internal xl/worksheets/sheet2.xml): add three new rows ---
$wsDropdowns = $wb.Worksheets.Item("Sheet4")
$wsDropdowns.Range("B24").Value = "Dropdowns"
$wsDropdowns.Range("B25").Value = "bank description"
$wsDropdowns.Range("B26").Value = "remove features to keep only 3"
[void]$wsDropdowns.Range("B25").Select()

# --- "Hours worked" tab (internal xl/worksheets/sheet4.xml): log new hours entry ---
$wsHours = $wb.Worksheets.Item("Hours worked")
$wsHours.Range("B15").Value = "9th Jan"
$wsHours.Range("C15").Value = 9
$wsHours.Range("D15").Value = "Updated dynamic dropdowns for:"
$wsHours.Range("D16").Value = "Refinance/buy new"
$wsHours.Range("D17").Value = "owner/investor"
$wsHours.Range("D18").Value = "principal&interest/interest"
$wsHours.Range("D19").Value = "offset/redraw"
$wsHours.Range("D20").Value = "updated javascript function to show hide dropdowns"

# Style the sub-bullet list (D16:D20) in italics, matching the new "note" style
$wsHours.Range("D16:D20").Font.Italic = $true

[void]$wsHours.Range("D21").Select()
